$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: populate C53:H53 (previously empty) with OHLCV + change values
$row53 = New-Object 'object[,]' 1,6
$row53[0,0] = 583.41
$row53[0,1] = 589.08000000000004
$row53[0,2] = 582.84
$row53[0,3] = 586.84
$row53[0,4] = 67808584
$row53[0,5] = 0.40680561042909069
$ws.Range("C53:H53").Value = $row53

# Row 54: populate A54:II54 (previously mostly empty) with a full new data row
$row54 = New-Object 'object[,]' 1,243
$row54[0,0] = 45790
$row54[0,1] = 45791
$row54[0,8] = 0.1822
$row54[0,9] = 3.8
$row54[0,10] = 600
$row54[0,11] = 328381800
$row54[0,12] = 30238
$row54[0,13] = 1572
$row54[0,14] = 31810
$row54[0,15] = 0.10983232514305166
$row54[0,16] = 0.023089327015383774
$row54[0,17] = 45793
$row54[0,18] = 0.17954227890308189
$row54[0,19] = 45807
$row54[0,20] = 0.067400951792040614
$row54[0,21] = 45828
$row54[0,22] = 0.36883820560211267
$row54[0,23] = 19.333333333333332
$row54[0,24] = 590
$row54[0,25] = 301461680
$row54[0,26] = -37674
$row54[0,27] = 28846
$row54[0,28] = 66520
$row54[0,29] = 0.10082847848428445
$row54[0,30] = 0.25904374886088255
$row54[0,31] = 45791
$row54[0,32] = 0.25904374886088255
$row54[0,33] = 45793
$row54[0,34] = 0.26675599597680588
$row54[0,35] = 45828
$row54[0,36] = 0.1155757015275987
$row54[0,37] = 14
$row54[0,38] = 595
$row54[0,39] = 135586220
$row54[0,40] = 13761
$row54[0,41] = 2017
$row54[0,42] = 15778
$row54[0,43] = 0.045348889006508082
$row54[0,44] = 0.15016989419597226
$row54[0,45] = 45791
$row54[0,46] = 0.15016989419597226
$row54[0,47] = 45793
$row54[0,48] = 0.22742538270591531
$row54[0,49] = 45828
$row54[0,50] = 0.18292533732186761
$row54[0,51] = 14
$row54[0,52] = 588
$row54[0,53] = 117031404
$row54[0,54] = 14260
$row54[0,55] = 16726
$row54[0,56] = 30986
$row54[0,57] = 0.039142946460722965
$row54[0,58] = 0.59904762989713034
$row54[0,59] = 45791
$row54[0,60] = 0.59904762989713034
$row54[0,61] = 45792
$row54[0,62] = 0.07602252575538368
$row54[0,63] = 45793
$row54[0,64] = 0.19416502432189198
$row54[0,65] = 2
$row54[0,66] = 589
$row54[0,67] = 95033383
$row54[0,68] = 15046
$row54[0,69] = 9570
$row54[0,70] = 24616
$row54[0,71] = 0.03178537123890593
$row54[0,72] = 0.53978074726917724
$row54[0,73] = 45791
$row54[0,74] = 0.53978074726917724
$row54[0,75] = 45793
$row54[0,76] = 0.16566417642516734
$row54[0,77] = 45800
$row54[0,78] = 0.10720661791104968
$row54[0,79] = 4.666666666666667
$row54[0,80] = 550
$row54[0,81] = -60958700
$row54[0,82] = 0.039383389241365682
$row54[0,83] = -1235
$row54[0,84] = 20092
$row54[0,85] = 21327
$row54[0,86] = 0
$row54[0,87] = 45828
$row54[0,88] = 0.64445947176092699
$row54[0,89] = 45838
$row54[0,90] = 0.033648906979466042
$row54[0,91] = 45856
$row54[0,92] = 0.2751667815389417
$row54[0,93] = 50.666666666666664
$row54[0,94] = 575
$row54[0,95] = -60818900
$row54[0,96] = 0.039293069109605276
$row54[0,97] = -7956
$row54[0,98] = 45011
$row54[0,99] = 52967
$row54[0,100] = 0.047708620184958313
$row54[0,101] = 45793
$row54[0,102] = 0.16074894023005953
$row54[0,103] = 45800
$row54[0,104] = 0.084282879972245414
$row54[0,105] = 45828
$row54[0,106] = 0.22904040677819096
$row54[0,107] = 17
$row54[0,108] = 565
$row54[0,109] = -37568545
$row54[0,110] = 0.024271787800047611
$row54[0,111] = -6511
$row54[0,112] = 1803
$row54[0,113] = 8314
$row54[0,114] = 0
$row54[0,115] = 45793
$row54[0,116] = 0.29348767142263471
$row54[0,117] = 45828
$row54[0,118] = 0.12971845130904558
$row54[0,119] = 45919
$row54[0,120] = 0.14776977420695386
$row54[0,121] = 56.666666666666664
$row54[0,122] = 555
$row54[0,123] = -32290455
$row54[0,124] = 0.020861789343371868
$row54[0,125] = -1032
$row54[0,126] = 57
$row54[0,127] = 1089
$row54[0,128] = 0
$row54[0,129] = 45800
$row54[0,130] = 0.09217407741908383
$row54[0,131] = 45828
$row54[0,132] = 0.35886664846359034
$row54[0,133] = 45856
$row54[0,134] = 0.29064711936317239
$row54[0,135] = 38
$row54[0,136] = 577
$row54[0,137] = -27308833
$row54[0,138] = 0.017643329004169251
$row54[0,139] = -4115
$row54[0,140] = 2028
$row54[0,141] = 6143
$row54[0,142] = 0.49560961566143658
$row54[0,143] = 45791
$row54[0,144] = 0.49560961566143658
$row54[0,145] = 45793
$row54[0,146] = 0.1628400748524543
$row54[0,147] = 45807
$row54[0,148] = 0.12073556931049374
$row54[0,149] = 7
$row54[0,150] = 590
$row54[0,151] = 397763840
$row54[0,152] = -37674
$row54[0,153] = 28846
$row54[0,154] = 66520
$row54[0,155] = 0.087658086184154688
$row54[0,156] = 349612760
$row54[0,157] = 0.11693334505895177
$row54[0,158] = 0.25904374886088255
$row54[0,159] = 45791
$row54[0,160] = 0.25904374886088255
$row54[0,161] = 45793
$row54[0,162] = 0.26675599597680588
$row54[0,163] = 45828
$row54[0,164] = 0.1155757015275987
$row54[0,165] = 14
$row54[0,166] = -48151080
$row54[0,167] = 0.031108811802616168
$row54[0,168] = 0.19071950203401461
$row54[0,169] = 45791
$row54[0,170] = 0.19071950203401461
$row54[0,171] = 45793
$row54[0,172] = 0.12957653286281429
$row54[0,173] = 45856
$row54[0,174] = 0.23013772484438563
$row54[0,175] = 23.333333333333332
$row54[0,176] = 600
$row54[0,177] = 374124600
$row54[0,178] = 30238
$row54[0,179] = 1572
$row54[0,180] = 31810
$row54[0,181] = 0.082448536373774944
$row54[0,182] = 351253200
$row54[0,183] = 0.11748201535510602
$row54[0,184] = 0.023089327015383774
$row54[0,185] = 45793
$row54[0,186] = 0.17954227890308189
$row54[0,187] = 45807
$row54[0,188] = 0.067400951792040614
$row54[0,189] = 45828
$row54[0,190] = 0.36883820560211267
$row54[0,191] = 19.333333333333332
$row54[0,192] = -22871400
$row54[0,193] = 0.014776451084012142
$row54[0,194] = 0.00078700910307195882
$row54[0,195] = 45884
$row54[0,196] = 0.22322201526797661
$row54[0,197] = 45919
$row54[0,198] = 0.19596526666491776
$row54[0,199] = 46038
$row54[0,200] = 0.19677850940475877
$row54[0,201] = 157
$row54[0,202] = 580
$row54[0,203] = 309104040
$row54[0,204] = -12595
$row54[0,205] = 41243
$row54[0,206] = 53838
$row54[0,207] = 0.068119486623495984
$row54[0,208] = 143945560
$row54[0,209] = 0.048144798368297673
$row54[0,210] = 0
$row54[0,211] = 45814
$row54[0,212] = 0.098202125859248451
$row54[0,213] = 45828
$row54[0,214] = 0.57440507369591676
$row54[0,215] = 45919
$row54[0,216] = 0.059673948956813952
$row54[0,217] = 63.666666666666664
$row54[0,218] = -165158480
$row54[0,219] = 0.10670340253896997
$row54[0,220] = 0.16788057143659835
$row54[0,221] = 45791
$row54[0,222] = 0.16788057143659835
$row54[0,223] = 45793
$row54[0,224] = 0.31072918568880026
$row54[0,225] = 45828
$row54[0,226] = 0.086927053336891935
$row54[0,227] = 14
$row54[0,228] = 0
$row54[0,229] = 81682
$row54[0,230] = 382803
$row54[0,231] = 2989846564.5
$row54[0,232] = -1547827680
$row54[0,233] = 1442018884.5
$row54[0,234] = 1.9316404552863404
$row54[0,235] = 4537674244.5
$row54[0,236] = 0.17378774797592239
$row54[0,237] = 45791
$row54[0,238] = 0.17378774797592239
$row54[0,239] = 45793
$row54[0,240] = 0.18487024867789625
$row54[0,241] = 45828
$row54[0,242] = 0.17050883366909791
$ws.Range("A54:II54").Value = $row54

# Update the active cell / selection in the bottom-right frozen pane
$ws.Range("F65").Select()
